$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1007786943"
$ws.Range("D16").Value = "ALFONSO GONZALEZ CONTRERAS"
$ws.Range("C17").Value = "1001898009"
$ws.Range("D17").Value = "ROBERTO VILLA PIMENTEL"
